# Add team record (Wins/Losses/Ties) columns to the roster sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) by
# copying the style from an existing header cell onto the new ones.
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Every player on this sheet shares the same 1996 team record: 83-79-0.
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 29).Value = 83  # AC - Wins
    $ws.Cells.Item($r, 30).Value = 79  # AD - Losses
    $ws.Cells.Item($r, 31).Value = 0   # AE - Ties
}
